# Scheduled-runner update: refreshes market-price / profit columns (H-N)
# for specific Leve rows across several crafting-class sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 817.5
$ws.Range("I17").Value = 175
$ws.Range("J17").Value = 844.2708
$ws.Range("K17").Value = 525
$ws.Range("L17").Value = 2532.8124
$ws.Range("M17").Value = -357
$ws.Range("N17").Value = -2868.8124

# Row 136: I Like Big Brush and I Cannot Lie / Dark Mahogany Round Brush
$ws.Range("H136").Value = 66806.875
$ws.Range("J136").Value = 66806.875
$ws.Range("L136").Value = 66806.875
$ws.Range("N136").Value = -77006.875

# Row 139: Something Salty and Ceremonial / Gomphotherium Codex
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 126717.5
$ws.Range("I134").Value = 159526.69
$ws.Range("J134").Value = 2042.6
$ws.Range("K134").Value = 478580.07
$ws.Range("L134").Value = 6127.799999999999
$ws.Range("M134").Value = -476045.07
$ws.Range("N134").Value = -11197.8

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2109.6428
$ws.Range("I31").Value = 1587.7368
$ws.Range("J31").Value = 3211.4443
$ws.Range("K31").Value = 1587.7368
$ws.Range("L31").Value = 3211.4443
$ws.Range("M31").Value = -1292.7368
$ws.Range("N31").Value = -3801.4443

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2109.6428
$ws.Range("I34").Value = 1587.7368
$ws.Range("J34").Value = 3211.4443
$ws.Range("K34").Value = 1587.7368
$ws.Range("L34").Value = 3211.4443
$ws.Range("M34").Value = -1385.7368
$ws.Range("N34").Value = -3615.4443

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 3553.0286
$ws.Range("I134").Value = 3853.5667
$ws.Range("J134").Value = 1749.8
$ws.Range("K134").Value = 11560.7001
$ws.Range("L134").Value = 5249.4
$ws.Range("M134").Value = -9025.7001
$ws.Range("N134").Value = -10319.4

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Range("H2").Value = 83461.336
$ws.Range("I2").Value = 111255.445
$ws.Range("J2").Value = 79
$ws.Range("K2").Value = 667532.67
$ws.Range("L2").Value = 474
$ws.Range("M2").Value = -667419.67
$ws.Range("N2").Value = -700

# Row 64: The Aroma of Faith / Baked Onion Soup
$ws.Range("H64").Value = 4629.4
$ws.Range("I64").Value = 683.2
$ws.Range("J64").Value = 6207.88
$ws.Range("K64").Value = 2049.6
$ws.Range("L64").Value = 18623.64
$ws.Range("M64").Value = -1779.6
$ws.Range("N64").Value = -19163.64

# Row 67: Soup's On (L) / Baked Onion Soup
$ws.Range("H67").Value = 4629.4
$ws.Range("I67").Value = 683.2
$ws.Range("J67").Value = 6207.88
$ws.Range("K67").Value = 2049.6
$ws.Range("L67").Value = 18623.64
$ws.Range("M67").Value = -1113.6
$ws.Range("N67").Value = -20495.64

# Row 70: Persona non Gratin / Dhalmel Gratin
$ws.Range("H70").Value = 6566.5713
$ws.Range("I70").Value = 5455.5
$ws.Range("J70").Value = 7576.636
$ws.Range("K70").Value = 16366.5
$ws.Range("L70").Value = 22729.908
$ws.Range("M70").Value = -16051.5
$ws.Range("N70").Value = -23359.908

# Row 73: Recipe for Disaster (L) / Dhalmel Gratin
$ws.Range("H73").Value = 6566.5713
$ws.Range("I73").Value = 5455.5
$ws.Range("J73").Value = 7576.636
$ws.Range("K73").Value = 16366.5
$ws.Range("L73").Value = 22729.908
$ws.Range("M73").Value = -15274.5
$ws.Range("N73").Value = -24913.908

# Row 82: Persuasion of a Higher Power / Baked Pipira Pira
$ws.Range("H82").Value = 7375.85
$ws.Range("I82").Value = 849.75
$ws.Range("J82").Value = 9007.375
$ws.Range("K82").Value = 2549.25
$ws.Range("L82").Value = 27022.125
$ws.Range("M82").Value = -2143.25
$ws.Range("N82").Value = -27834.125

# Row 85: Loaves and Fishes (L) / Baked Pipira Pira
$ws.Range("H85").Value = 7375.85
$ws.Range("I85").Value = 849.75
$ws.Range("J85").Value = 9007.375
$ws.Range("K85").Value = 2549.25
$ws.Range("L85").Value = 27022.125
$ws.Range("M85").Value = -1145.25
$ws.Range("N85").Value = -29830.125

# Row 88: Don't Let It Fall Apart / Liver-cheese Sandwich
$ws.Range("H88").Value = 3436.3635
$ws.Range("J88").Value = 3436.3635
$ws.Range("L88").Value = 10309.0905
$ws.Range("N88").Value = -11165.0905

# Row 91: Better Come Back with a Sandwich (L) / Liver-cheese Sandwich
$ws.Range("H91").Value = 3436.3635
$ws.Range("J91").Value = 3436.3635
$ws.Range("L91").Value = 10309.0905
$ws.Range("N91").Value = -13273.0905

# Row 94: All You Can Stomach / Baklava
$ws.Range("H94").Value = 8545.777
$ws.Range("I94").Value = 5512
$ws.Range("J94").Value = 9412.571
$ws.Range("K94").Value = 16536
$ws.Range("L94").Value = 28237.713
$ws.Range("M94").Value = -15860
$ws.Range("N94").Value = -29589.713

# Row 100: Souper / Gameni
$ws.Range("H100").Value = 11912578
$ws.Range("J100").Value = 11912578
$ws.Range("L100").Value = 35737734
$ws.Range("N100").Value = -35739356

# Row 103: West Meats East / Nomad Meat Pie
$ws.Range("H103").Value = 1421.4
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 1469.8948
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 4409.6844
$ws.Range("M103").Value = -621
$ws.Range("N103").Value = -6167.6844

# Row 106: Herky Jerky / Jerked Jhammel
$ws.Range("H106").Value = 3658.3333
$ws.Range("J106").Value = 3658.3333
$ws.Range("L106").Value = 10974.9999
$ws.Range("N106").Value = -12866.9999

# Row 112: Sweet Tooth / Caramels
$ws.Range("H112").Value = 46708430
$ws.Range("I112").Value = 1184.6666
$ws.Range("J112").Value = 53077600
$ws.Range("K112").Value = 3553.9998
$ws.Range("L112").Value = 159232800
$ws.Range("M112").Value = -2445.9998
$ws.Range("N112").Value = -159235016

# Row 120: A Happy End / Paella
$ws.Range("H120").Value = 14586.429
$ws.Range("I120").Value = 4752.5
$ws.Range("J120").Value = 18520
$ws.Range("K120").Value = 14257.5
$ws.Range("L120").Value = 55560
$ws.Range("M120").Value = -9419.5
$ws.Range("N120").Value = -65236

# Row 125: At Any Temperature / Borscht
$ws.Range("H125").Value = 2641.7646
$ws.Range("I125").Value = 1140
$ws.Range("J125").Value = 2787.0967
$ws.Range("K125").Value = 3420
$ws.Range("L125").Value = 8361.2901
$ws.Range("M125").Value = 1500
$ws.Range("N125").Value = -18201.2901

$ws = $wb.Worksheets.Item("LTW")
# Row 56: Hold On Tight / Raptorskin Smithy's Gloves
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1702.1333
$ws.Range("I61").Value = 1425.5385
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 1425.5385
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -1223.5385
$ws.Range("N61").Value = -3904

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1702.1333
$ws.Range("I113").Value = 1425.5385
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1425.5385
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 744.4614999999999
$ws.Range("N113").Value = -7840

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3197.5
$ws.Range("I122").Value = 2393.3333
$ws.Range("J122").Value = 3855.4546
$ws.Range("K122").Value = 7179.999899999999
$ws.Range("L122").Value = 11566.3638
$ws.Range("M122").Value = -4729.999899999999
$ws.Range("N122").Value = -16466.3638

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 1708.1608
$ws.Range("I136").Value = 1457.4375
$ws.Range("J136").Value = 3212.5
$ws.Range("K136").Value = 4372.3125
$ws.Range("L136").Value = 9637.5
$ws.Range("M136").Value = -1822.3125
$ws.Range("N136").Value = -14737.5

$ws = $wb.Worksheets.Item("WVR")
# Row 58: Seeing It Through to the End / Woolen Smock
$ws.Range("H58").Value = 9000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9692

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 46586.184
$ws.Range("I113").Value = 71778.28999999999
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 215334.87
$ws.Range("L113").Value = 7500
$ws.Range("M113").Value = -213164.87
$ws.Range("N113").Value = -11840

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1567.5209
$ws.Range("I132").Value = 1446.262
$ws.Range("J132").Value = 2416.3333
$ws.Range("K132").Value = 4338.786
$ws.Range("L132").Value = 7248.999899999999
$ws.Range("M132").Value = -1808.786
$ws.Range("N132").Value = -12308.9999
